# "New" button on the settings sheet: add a fresh worksheet (Excel's
# default "Sheet1") after the existing sheets and hard-code the default
# population-stats values/formulas that the button seeds it with.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

# row 1 - starting/selected population
$ws.Range("A1").Value = 153
$ws.Range("B1").Value = "selected"

# row 2 - population that came in from the right, plus the running ratio
$ws.Range("A2").Value = 162
$ws.Range("B2").Value = "popstats"
$ws.Range("C2").Formula = "=A2/A3"

# row 3 - sum of rows 1 & 2
$ws.Range("A3").Formula = "=A1+A2"
$ws.Range("B3").Value = "sum"

# row 4 - popright total
$ws.Range("A4").Value = 366
$ws.Range("B4").Value = "popright"

# row 5 - filler, difference between popright total and the sum
$ws.Range("A5").Formula = "=A4-A3"
$ws.Range("B5").Value = "filler"

# row 7 - literal (quote-prefixed) text note of the formula being modeled
$ws.Range("A7").Value = "'=popright-filler(153/315"

$ws.Activate()
$ws.Range("C3").Select()
